# Update odds data for rows 2 and 7 on the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 2.57
$ws.Range("I2").Value = 2.67
$ws.Range("J2").Value = 3.15
$ws.Range("L2").Value = 3.3
$ws.Range("M2").Value = 9
$ws.Range("O2").Value = 1.28
$ws.Range("Q2").Value = 1.9
$ws.Range("R2").Value = 1.82
$ws.Range("U2").Value = 1.62
$ws.Range("V2").Value = 2.02
$ws.Range("W2").Value = 8.75
$ws.Range("X2").Value = 13.5
$ws.Range("Y2").Value = 9.5
$ws.Range("Z2").Value = 30
$ws.Range("AA2").Value = 21
$ws.Range("AB2").Value = 28
$ws.Range("AD2").Value = 6
$ws.Range("AE2").Value = 12.5
$ws.Range("AH2").Value = 9
$ws.Range("AJ2").Value = 9.75
$ws.Range("AL2").Value = 22
$ws.Range("AM2").Value = 29
$ws.Range("AN2").Value = 4.45
$ws.Range("AO2").Value = 14
$ws.Range("AQ2").Value = 60
$ws.Range("AU2").Value = 6.7
$ws.Range("AW2").Value = 4.6
$ws.Range("AX2").Value = 14.5
$ws.Range("AY2").Value = 22
$ws.Range("AZ2").Value = 65
$ws.Range("BA2").Value = 100

# Row 7 updates
$ws.Range("G7").Value = 3.25
$ws.Range("H7").Value = 3
$ws.Range("I7").Value = 2.2
$ws.Range("N7").Value = 8
$ws.Range("U7").Value = 1.95
$ws.Range("V7").Value = 1.8
$ws.Range("Y7").Value = 13
$ws.Range("Z7").Value = 34
$ws.Range("AI7").Value = 10
$ws.Range("AO7").Value = 19
$ws.Range("AP7").Value = 29
